$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "42.958.33"
Set-TextValue "E2" "  -0.25%  "
Set-TextValue "D3" "2.363.35"
Set-TextValue "E3" "  +1.29%  "
Set-TextValue "D5" "303.25"
Set-TextValue "E5" "  +0.06%  "
Set-TextValue "E6" "  -0.68%  "
Set-TextValue "E7" "  -0.06%  "
Set-TextValue "E8" "  -0.72%  "
Set-TextValue "E9" "  -2.70%  "
Set-TextValue "D10" "34.17"
Set-TextValue "E10" "  -0.87%  "
Set-TextValue "E11" "  +2.42%  "
Set-TextValue "E12" "  -0.12%  "
Set-TextValue "D13" "18.44"
Set-TextValue "E13" "  -3.69%  "
Set-TextValue "D14" "6.72"
Set-TextValue "E14" "  -0.58%  "
Set-TextValue "D15" "2.728.12"
Set-TextValue "E15" "  +1.22%  "
Set-TextValue "D16" "2.366.19"
Set-TextValue "E16" "  +0.72%  "
Set-TextValue "D17" "0.792"
Set-TextValue "E17" "  -0.16%  "
Set-TextValue "D18" "42.950.27"
Set-TextValue "E18" "  -0.15%  "
Set-TextValue "D19" "11.91"
Set-TextValue "E19" "  -2.67%  "
Set-TextValue "E20" "  +1.02%  "
Set-TextValue "D21" "0.0₃0886"
Set-TextValue "E21" "  -0.85%  "
Set-TextValue "D22" "68.09"
Set-TextValue "E22" "  +0.09%  "
Set-TextValue "D23" "235.01"
Set-TextValue "E23" "  -0.79%  "
Set-TextValue "D24" "2.17"
Set-TextValue "D25" "2.44"
Set-TextValue "E25" "  +0.18%  "
Set-TextValue "E26" "  -0.07%  "
Set-TextValue "D27" "24.55"
Set-TextValue "E27" "  -0.91%  "
Set-TextValue "E28" "  +15.16%  "
Set-TextValue "D29" "9.34"
Set-TextValue "E29" "  +1.87%  "
Set-TextValue "D30" "32.18"
Set-TextValue "E30" "  -0.44%  "
Set-TextValue "E31" "  +0.01%  "
Set-TextValue "E32" "  -0.59%  "
Set-TextValue "D33" "17.50"
Set-TextValue "E33" "  -2.57%  "
Set-TextValue "D34" "0.0714"
Set-TextValue "E34" "  +1.44%  "
Set-TextValue "E35" "  +2.88%  "
Set-TextValue "D36" "1.84"
Set-TextValue "E36" "  +1.01%  "
Set-TextValue "D37" "126.60"
Set-TextValue "E37" "  -11.23%  "
Set-TextValue "D38" "4.32"
Set-TextValue "E38" "  -2.76%  "
Set-TextValue "D39" "2.84"
Set-TextValue "E39" "  +3.07%  "
Set-TextValue "D40" "2.26"
Set-TextValue "E40" "  -1.70%  "
Set-TextValue "E41" "  -1.37%  "
Set-TextValue "D42" "21.42"
Set-TextValue "E42" "  -4.30%  "
Set-TextValue "D43" "1.931.97"
Set-TextValue "E43" "  -0.21%  "
Set-TextValue "E44" "  -0.37%  "
Set-TextValue "E45" "  +3.91%  "
Set-TextValue "D46" "9.24"
Set-TextValue "E46" "  -8.78%  "
Set-TextValue "D47" "2.71"
Set-TextValue "E47" "  -2.24%  "
Set-TextValue "D48" "2.588.98"
Set-TextValue "E48" "  +0.99%  "
Set-TextValue "E49" "  +1.22%  "
Set-TextValue "D50" "71.49"
Set-TextValue "E50" "  -2.60%  "
Set-TextValue "E51" "  +0.84%  "
